$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Bring the classic Rexx interface rows in sync with the object-oriented
# interface: several functions that were only documented as "possible
# future enhancement" / "used to create new SQL functions..." are now
# implemented (both oo and classic), and a couple of comments are
# corrected / reworded.
# ---------------------------------------------------------------------------

# Donor row (100) already carries the plain "implemented oo and classic"
# look: A/B styled, C unstyled - copy that formatting onto the rows whose
# status changed to "implemented".
$donor = $ws.Range("A100:C100")

# Row 65 - sqlite3_context_db_handle: now implemented oo and classic
$donor.Copy()
$ws.Range("A65:C65").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B65").Value = " implemented oo and classic"
$ws.Range("C65").Value = "Implemented indirectly, framework passes db to UDFs"

# Row 5 - sqlite3_aggregate_context: now implemented oo and classic
$donor.Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B5").Value = " implemented oo and classic"
$ws.Range("C5").Value = "Implemented indirectly, framework uses in implementation of UDFs"

# Row 183 - sqlite3_user_data: now implemented oo and classic
$donor.Copy()
$ws.Range("A183:C183").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B183").Value = " implemented oo and classic"
$ws.Range("C183").Value = "Implemented indirectly, used by implementation code"

# Row 36 - sqlite3_collation_needed: classic interface caught up too,
# and the function name cell drops its "NEED" highlight style.
$srcA = $ws.Range("A100")
$srcA.Copy()
$ws.Range("A36").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B36").Value = " implemented oo and classic"

# Row 94 - sqlite3_free: B94 picks up the standard "implemented" styling
$ws.Range("B94:B94").Copy() | Out-Null
$srcB = $ws.Range("B100")
$srcB.Copy()
$ws.Range("B94").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 112 - sqlite3_mprintf: reworded comment
$ws.Range("C112").Value = "Implemented indirectly in the equote() method and ooSQLiteEnquote()"

# ---------------------------------------------------------------------------
# Cosmetic worksheet changes: widen comments column, move the frozen-pane
# scroll position / active selection.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 68.140625

$ws.Range("A145").Select()
$excel.ActiveWindow.ScrollRow = 56
